$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.49019999999999
$ws.Range("D3").Value = -5.832099999999998
$ws.Range("D5").Value = -7.989299999999994
$ws.Range("B9").Value = 8.585500000000007
$ws.Range("D11").Value = -8.223700000000003
$ws.Range("D12").Value = -8.4496
$ws.Range("B13").Value = 5.354200000000002
$ws.Range("B16").Value = 9.015700000000006
$ws.Range("B18").Value = 4.986200000000002
$ws.Range("B20").Value = 5.7318
$ws.Range("D21").Value = -7.413299999999995
